$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tuesday")

# Tuesday speeches pulled from Cision (added below the existing rows)
$rows = @(
    @("Day",        "Sharon",    "rnc.day.txt",          "Tuesday", "speech", "Cision"),
    @("White",       "Dana",      "rnc.white.txt",        "Tuesday", "speech", "Cision"),
    @("Hutchinson",  "Asa",       "rnc.hutchinson.txt",   "Tuesday", "speech", "Cision"),
    @("Rutledge",    "Leslie",    "rnc.rutledge.txt",     "Tuesday", "speech", "Cision"),
    @("Mukasey",     "Michael",   "rnc.mukasey.txt",      "Tuesday", "speech", "Cision"),
    @("Wist",        "Andy",      "rnc.wist.txt",         "Tuesday", "speech", "Cision"),
    @("Johnson",     "Ron",       "rnc.johnson.txt",      "Tuesday", "speech", "Cision"),
    @("Cox",         "Chris",     "rnc.cox.txt",          "Tuesday", "speech", "Cision"),
    @("Gulbis",      "natalie",   "rnc.gulbis.txt",       "Tuesday", "speech", "Cision"),
    @("Sullivan",    "Dana",      "rnc.sullivandan.txt",  "Tuesday", "speech", "Cision"),
    @("Woolard",     "Kerry",     "rnc.woolard.txt",      "Tuesday", "speech", "Cision"),
    @("Capito",      "Shelly",    "rnc.capito.txt",       "Tuesday", "speech", "Cision"),
    @("Brown",       "Kimberlin", "rnc.brownkim.txt",     "Tuesday", "speech", "Cision")
)

$startRow = 9
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

$ws.Range("F22").Select()
